$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename the "Requested quantity" header on the existing sheets.
# ---------------------------------------------------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# ---------------------------------------------------------------------------
# 2. Add the new "PO Forecast" sheet at the end of the workbook.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Header row values
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Match header formatting (bold, border, centered/top aligned) used on the
# other sheets by copying the existing header cell's format.
$wsWeekly.Range("B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Forecast data rows
$data = @(
    @(45319.99999999999, 138, 63.62307876117139, 220.5536648455232),
    @(45326.99999999999, 150, 68.66558983142308, 229.6968448811789),
    @(45333.99999999999, 162, 80.22019230129784, 243.9250967656285),
    @(45340.99999999999, 174, 90.83560494714649, 252.0580868812682),
    @(45347.99999999999, 185, 100.8617333238307, 264.5356699333664),
    @(45354.99999999999, 197, 111.8195903498278, 275.7146433330383),
    @(45361.99999999999, 209, 129.9936466623081, 288.3293334563759),
    @(45368.99999999999, 221, 143.036447283026, 306.56807438889),
    @(45375.99999999999, 233, 155.7997291019036, 323.6878979806847),
    @(45382.99999999999, 245, 161.6209536084453, 324.3971979621673),
    @(45389.99999999999, 257, 171.2022205292101, 338.1707503161567)
)

$row = 2
foreach ($r in $data) {
    $wsForecast.Cells.Item($row, 1).Value = $r[0]
    $wsForecast.Cells.Item($row, 2).Value = $r[1]
    $wsForecast.Cells.Item($row, 3).Value = $r[2]
    $wsForecast.Cells.Item($row, 4).Value = $r[3]
    $row++
}

# Match the date-column formatting (column A, rows 2-12) used on the other
# sheets by copying the existing date cell's format.
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A12").PasteSpecial(-4122)
$excel.CutCopyMode = $false
